# Weekly cryptos-list refresh (GitHub Actions bot).
# Updates Price (D) / Volume(1h) (E) text cells with freshly scraped values,
# and refreshes the Mantle / InjectiveProtocol ranking swap at rows 43-44.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume cells are plain text in the sheet (e.g. "63.467.93", "0.581",
# "  -1.00%  "). Values that look like a plain decimal number would be auto-
# converted to a float by Excel and lose trailing zeros (e.g. '5.40' -> 5.4),
# so those are written with a leading apostrophe to force text, exactly as a
# person retyping the figure in Excel would do.

$ws.Range("D2").Value = '63.513.64'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").Value = '2.689.39'
$ws.Range("E3").Value = '  -2.50%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''553.96'
$ws.Range("E5").Value = '  -3.82%  '

$ws.Range("D6").Value = '''158.12'
$ws.Range("E6").Value = '  -0.77%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '''0.582'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("E9").Value = '  -4.11%  '

$ws.Range("E11").Value = '  -4.33%  '

$ws.Range("D12").Value = '''5.40'
$ws.Range("E12").Value = '  -8.31%  '

$ws.Range("D13").Value = '3.166.16'
$ws.Range("E13").Value = '  -2.49%  '

$ws.Range("D14").Value = '''26.40'
$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("D15").Value = '63.382.26'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("E16").Value = '  -4.23%  '

$ws.Range("D17").Value = '2.692.59'
$ws.Range("E17").Value = '  -2.57%  '

$ws.Range("D18").Value = '''12.09'
$ws.Range("E18").Value = '  -0.57%  '

$ws.Range("E19").Value = '  -4.85%  '

$ws.Range("D20").Value = '''343.39'
$ws.Range("E20").Value = '  -4.54%  '

$ws.Range("E21").Value = '  -4.83%  '

$ws.Range("D22").Value = '''0.996'
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").Value = '''0.507'
$ws.Range("E23").Value = '  -3.96%  '

$ws.Range("D24").Value = '''63.97'
$ws.Range("E24").Value = '  -1.69%  '

$ws.Range("D25").Value = '''0.170'
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").Value = '''8.18'
$ws.Range("E27").Value = '  -4.34%  '

$ws.Range("D28").Value = '0.0₃0857'
$ws.Range("E28").Value = '  -5.26%  '

$ws.Range("E29").Value = '  -0.66%  '

$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").Value = '''7.02'
$ws.Range("E31").Value = '  -4.81%  '

$ws.Range("D32").Value = '''165.45'
$ws.Range("E32").Value = '  -2.64%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").Value = '''4.79'
$ws.Range("E34").Value = '  -3.18%  '

$ws.Range("D35").Value = '''19.58'
$ws.Range("E35").Value = '  -3.26%  '

$ws.Range("E36").Value = '  -3.75%  '

$ws.Range("E37").Value = '  -1.95%  '

$ws.Range("D38").Value = '''340.06'
$ws.Range("E38").Value = '  -2.38%  '

$ws.Range("D39").Value = '''0.949'
$ws.Range("E39").Value = '  -5.86%  '

$ws.Range("D40").Value = '''6.07'
$ws.Range("E40").Value = '  -3.90%  '

$ws.Range("D41").Value = '''38.18'
$ws.Range("E41").Value = '  -2.33%  '

$ws.Range("E42").Value = '  -6.03%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.623'
$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''20.82'
$ws.Range("E44").Value = '  -5.51%  '

$ws.Range("E45").Value = '  -5.94%  '

$ws.Range("D46").Value = '''0.0564'
$ws.Range("E46").Value = '  -4.06%  '

$ws.Range("D47").Value = '''0.999'
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("D49").Value = '''130.22'
$ws.Range("E49").Value = '  -5.19%  '

$ws.Range("E50").Value = '  -3.85%  '

$ws.Range("D51").Value = '2.096.98'
$ws.Range("E51").Value = '  -1.58%  '
